# Updated RAD Test Scripts and Test Data for Existing Liability.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Existing" sheet (sheet2): rename the "Existing Liability" label
#    used in column D (rows 2-19) to the new, longer wording, and widen
#    column D so the new text fits.
# ---------------------------------------------------------------------
$wsExisting = $wb.Worksheets.Item("Existing")
$newLabel = "Existing Liability with Notice/Invoice Number"
for ($r = 2; $r -le 19; $r++) {
    $wsExisting.Range("D" + $r).Value = $newLabel
}
$wsExisting.Range("D1").ColumnWidth = 51.1666666667

# ---------------------------------------------------------------------
# 2. "Personal_IND" sheet (sheet6): fill in the missing "Y" markers in
#    column C for every row (row 3 already had one).
# ---------------------------------------------------------------------
$wsInd = $wb.Worksheets.Item("Personal_IND")
foreach ($r in 2,4,5,6,7,8,9) {
    $cell = $wsInd.Range("C" + $r)
    $cell.Value = "Y"
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# 3. "Personal_JNT" sheet (sheet7): fill in the missing "Y" markers in
#    column C for every row (row 3 already had one).
# ---------------------------------------------------------------------
$wsJnt = $wb.Worksheets.Item("Personal_JNT")
foreach ($r in 2,4,5,6) {
    $cell = $wsJnt.Range("C" + $r)
    $cell.Value = "Y"
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# 4. "Personal_EL" sheet (sheet8): rename the "Existing Liability" label
#    used in column D (rows 2-3) to match the new wording used above.
# ---------------------------------------------------------------------
$wsEl = $wb.Worksheets.Item("Personal_EL")
$wsEl.Range("D2").Value = $newLabel
$wsEl.Range("D3").Value = $newLabel

# ---------------------------------------------------------------------
# 5. Restore / update the various sheet selections that Excel records
#    as part of normal interactive editing.
# ---------------------------------------------------------------------
$wsExisting.Activate() | Out-Null
$wsExisting.Range("D2:D19").Select() | Out-Null

$wsInd.Activate() | Out-Null
$wsInd.Range("C2:C9").Select() | Out-Null

$wsJnt.Activate() | Out-Null
$wsJnt.Range("C2:C6").Select() | Out-Null

$wsNewTaxReturn = $wb.Worksheets.Item("NewTaxReturn")
$wsNewTaxReturn.Activate() | Out-Null
$wsNewTaxReturn.Range("C2:C52").Select() | Out-Null

# "Personal_EL" ends up as the active sheet/tab after the edit session.
$wsEl.Activate() | Out-Null
$wsEl.Range("D6").Select() | Out-Null
